# Update of 2025 data and RF changes
# Column I ("RF") values for rows 31-50 change from 7.030625 to 6.777941176470589

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRF = 6.777941176470589

for ($row = 31; $row -le 50; $row++) {
    $ws.Cells.Item($row, 9).Value = $newRF
}
